# Inserts a new weekly price record as row 80 in the "Vega Modelo de Temuco -
# Bruselas (repollito)" sheet, pushing the existing rows 80-135 down to 81-136.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 80; this shifts rows 80-135 -> 81-136
# and copies formatting (incl. the date number format on column D) from the
# row above, same as Excel's normal "Insert Row" behaviour.
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new data record.
$ws.Range("A80").Value2 = 10
$ws.Range("B80").Value2 = "Vega Modelo de Temuco"
$ws.Range("C80").Value2 = "La Araucanía"
$ws.Range("D80").Value2 = 44824
$ws.Range("E80").Value2 = 9
$ws.Range("F80").Value2 = 100112035
$ws.Range("G80").Value2 = "Bruselas (repollito)"
$ws.Range("H80").Value2 = "Sin especificar"
$ws.Range("I80").Value2 = "Primera"
$ws.Range("J80").Value2 = 35
$ws.Range("K80").Value2 = 24000
$ws.Range("L80").Value2 = 24000
$ws.Range("M80").Value2 = 24000
$ws.Range("N80").Value2 = "$/malla 10 kilos"
$ws.Range("O80").Value2 = "Provincia de Quillota"
$ws.Range("P80").Value2 = 2400
$ws.Range("Q80").Value2 = 10
$ws.Range("R80").Value2 = "Hortaliza"
